$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster values (player, positions, team) for rows 2..19.
$data = @(
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Zach Edey", "C", "Memphis Grizzlies"),
    @("Goga Bitadze", "C", "Orlando Magic"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Amen Thompson", "SG,SF,PF", "Houston Rockets"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
